$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4: "123.456" -> "abc" (keeps its existing quote-prefix cell style)
$ws.Range("B4").Value = "abc"
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats - re-apply quote-prefix style lost by the value write

# B5: "123,456" -> empty, but the leading-quote cell style must remain
$ws.Range("B5").Value = ""
$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122)  # xlPasteFormats

# F12: "pippo" -> empty, but the leading-quote cell style must remain
$ws.Range("F12").Value = ""
$ws.Range("B3").Copy()
$ws.Range("F12").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Move the active selection to C3
$ws.Range("C3").Select()

$wb.Save()
